$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "RLC_TestSuite_15_coords.json"
$ws.Range("C2").Value = 0.6683718037592058
$ws.Range("D2").Value = 0.03288079560865044
$ws.Range("E2").Value = 246.6090285508171
$ws.Range("F2").Value = 104.3507244353449
$ws.Range("G2").Value = 6.934318167201073
$ws.Range("H2").Value = 1.18310657981665
$ws.Range("I2").Value = 98.32916168066278
$ws.Range("J2").Value = 8.485788589154039
$ws.Range("K2").Value = 1974.6271835113
$ws.Range("L2").Value = 447.854291077642
$ws.Range("B3").Value = "RLC_TestSuite_5_coords.json"
$ws.Range("C3").Value = 0.3931285850552266
$ws.Range("D3").Value = 0.03933038252520592
$ws.Range("E3").Value = 283.665561285869
$ws.Range("F3").Value = 79.89337494428813
$ws.Range("G3").Value = 3.962233930035661
$ws.Range("H3").Value = 0.9068341014186596
$ws.Range("I3").Value = 53.1280445697862
$ws.Range("J3").Value = 7.584197045143473
$ws.Range("K3").Value = 1814.024412895181
$ws.Range("L3").Value = 541.8071399581345
$ws.Range("B4").Value = "RLC_TestSuite_6_coords.json"
$ws.Range("C4").Value = 0.7301567056761517
$ws.Range("D4").Value = 0.03796968373175965
$ws.Range("E4").Value = 254.7724547092789
$ws.Range("F4").Value = 101.5862304253018
$ws.Range("G4").Value = 8.480522323243946
$ws.Range("H4").Value = 1.094424690773556
$ws.Range("I4").Value = 79.89276019576512
$ws.Range("J4").Value = 5.586110148184757
$ws.Range("K4").Value = 1850.592509979493
$ws.Range("L4").Value = 551.0379645283992
$ws.Range("B5").Value = "RLC_TestSuite_16_coords.json"
$ws.Range("C5").Value = 0.65755222307837
$ws.Range("D5").Value = 0.03329008083799603
$ws.Range("E5").Value = 307.6824223778115
$ws.Range("F5").Value = 36.05022801550312
$ws.Range("G5").Value = 6.3701263573774
$ws.Range("H5").Value = 1.131667226335917
$ws.Range("I5").Value = 95.651036315366
$ws.Range("J5").Value = 8.881763762976199
$ws.Range("K5").Value = 2008.738676920388
$ws.Range("L5").Value = 455.0839777070725
$ws.Range("B6").Value = "RLC_TestSuite_13_coords.json"
$ws.Range("C6").Value = 0.6693760267210556
$ws.Range("D6").Value = 0.03843684480498205
$ws.Range("E6").Value = 304.238566197809
$ws.Range("F6").Value = 45.95581259151502
$ws.Range("G6").Value = 6.851210799073301
$ws.Range("H6").Value = 1.245929603609504
$ws.Range("I6").Value = 93.11036654287889
$ws.Range("J6").Value = 8.670055411828127
$ws.Range("K6").Value = 1964.245983783335
$ws.Range("L6").Value = 455.2759368087738
$ws.Range("B7").Value = "RLC_TestSuite_3_coords.json"
$ws.Range("C7").Value = 0.404240199470016
$ws.Range("D7").Value = 0.04003296177616816
$ws.Range("E7").Value = 292.3484603462716
$ws.Range("F7").Value = 69.46119984859916
$ws.Range("G7").Value = 4.267263310495603
$ws.Range("H7").Value = 1.130638279165224
$ws.Range("I7").Value = 61.8824364119553
$ws.Range("J7").Value = 10.81251953924441
$ws.Range("K7").Value = 1807.645961080228
$ws.Range("L7").Value = 411.0887461937414
$ws.Range("B8").Value = "RLC_TestSuite_4_coords.json"
$ws.Range("C8").Value = 0.3937393250207335
$ws.Range("D8").Value = 0.02680134075726888
$ws.Range("E8").Value = 313.7857471599922
$ws.Range("F8").Value = 3.610689132743742
$ws.Range("G8").Value = 3.725824241955376
$ws.Range("H8").Value = 0.978435762548258
$ws.Range("I8").Value = 58.27079011521325
$ws.Range("J8").Value = 8.995885128995349
$ws.Range("K8").Value = 1816.690542899132
$ws.Range("L8").Value = 525.393213605718
$ws.Range("B9").Value = "RLC_TestSuite_14_coords.json"
$ws.Range("C9").Value = 0.6575195364293721
$ws.Range("D9").Value = 0.04000077297106158
$ws.Range("E9").Value = 313.7857471599922
$ws.Range("F9").Value = 3.610689132743743
$ws.Range("G9").Value = 6.463571505040402
$ws.Range("H9").Value = 1.147839937829721
$ws.Range("I9").Value = 90.3457502593194
$ws.Range("J9").Value = 9.600147766892855
$ws.Range("K9").Value = 2002.585947699395
$ws.Range("L9").Value = 464.8120958379947
$ws.Range("B10").Value = "RLC_TestSuite_8_coords.json"
$ws.Range("C10").Value = 0.6869305662193269
$ws.Range("D10").Value = 0.03884275883781806
$ws.Range("E10").Value = 293.3331241396158
$ws.Range("F10").Value = 66.85212650105903
$ws.Range("G10").Value = 7.709683208373604
$ws.Range("H10").Value = 0.9614151099865829
$ws.Range("I10").Value = 83.79952503427606
$ws.Range("J10").Value = 8.037653852467582
$ws.Range("K10").Value = 1942.961703727733
$ws.Range("L10").Value = 530.2377417271765
$ws.Range("B11").Value = "RLC_TestSuite_11_coords.json"
$ws.Range("C11").Value = 0.6661033843559852
$ws.Range("D11").Value = 0.04595395665896158
$ws.Range("E11").Value = 300.828548226054
$ws.Range("F11").Value = 53.10417837563701
$ws.Range("G11").Value = 6.879017028017589
$ws.Range("H11").Value = 1.271965106170889
$ws.Range("I11").Value = 85.5901933192828
$ws.Range("J11").Value = 9.361992536648961
$ws.Range("K11").Value = 2061.936250000983
$ws.Range("L11").Value = 491.90475189757
$ws.Range("B12").Value = "RLC_TestSuite_1_coords.json"
$ws.Range("C12").Value = 0.4036031430896432
$ws.Range("D12").Value = 0.03539408072544233
$ws.Range("E12").Value = 230.0150402972412
$ws.Range("F12").Value = 113.311771562717
$ws.Range("G12").Value = 4.273633587705902
$ws.Range("H12").Value = 1.042881361262353
$ws.Range("I12").Value = 54.80078055408472
$ws.Range("J12").Value = 10.36305049668642
$ws.Range("K12").Value = 1764.250802200371
$ws.Range("L12").Value = 516.3542763928889
$ws.Range("B13").Value = "RLC_TestSuite_2_coords.json"
$ws.Range("C13").Value = 0.4196212419310727
$ws.Range("D13").Value = 0.02836546144492515
$ws.Range("E13").Value = 295.3580128891678
$ws.Range("F13").Value = 61.7653005216096
$ws.Range("G13").Value = 4.756207413689229
$ws.Range("H13").Value = 0.8529711366137388
$ws.Range("I13").Value = 50.92195218270336
$ws.Range("J13").Value = 7.476336486966217
$ws.Range("K13").Value = 1830.340359177729
$ws.Range("L13").Value = 585.7373403722705
$ws.Range("B14").Value = "RLC_TestSuite_12_coords.json"
$ws.Range("C14").Value = 0.6789749283276268
$ws.Range("D14").Value = 0.04480296539380799
$ws.Range("E14").Value = 261.9219457639265
$ws.Range("F14").Value = 95.42176020829332
$ws.Range("G14").Value = 7.225371297754854
$ws.Range("H14").Value = 1.264274758782104
$ws.Range("I14").Value = 95.27001131341763
$ws.Range("J14").Value = 9.319955523342189
$ws.Range("K14").Value = 1959.740084060797
$ws.Range("L14").Value = 382.4327001994033
$ws.Range("B15").Value = "RLC_TestSuite_7_coords.json"
$ws.Range("C15").Value = 0.7059342377626692
$ws.Range("D15").Value = 0.04353920837381556
$ws.Range("E15").Value = 250.7931071417807
$ws.Range("F15").Value = 105.4667130685086
$ws.Range("G15").Value = 8.072713404424013
$ws.Range("H15").Value = 1.312295334816975
$ws.Range("I15").Value = 85.35090447348904
$ws.Range("J15").Value = 8.62723206072824
$ws.Range("K15").Value = 1954.892265060154
$ws.Range("L15").Value = 485.6724998482531
$ws.Rows("16:19").Delete()
